$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.04537966666667
$ws.Range("H2").Value = 54.13613900000001
$ws.Range("I2").Value = 0.6797959733292525
$ws.Range("J2").Value = 0.6797959733292525
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.701354
$ws.Range("N2").Value = 8.104061999999999
$ws.Range("O2").Value = 0.02221077311549548
$ws.Range("P2").Value = 0.02221077311549548
$ws.Range("Q2").Value = 48.74695854406867
$ws.Range("R2").Value = 438.722626896618
$ws.Range("S2").Value = 0.01509879412844344
$ws.Range("T2").Value = 0.01509879412844344
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.04537966666667
$ws.Range("H3").Value = 54.13613900000001
$ws.Range("I3").Value = 0.6797959733292525
$ws.Range("J3").Value = 0.6797959733292525
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("N3").Value = 264.449356
$ws.Range("O3").Value = 0.7247753838328104
$ws.Range("P3").Value = 0.7247753838328105
$ws.Range("Q3").Value = 1590.696343875165
$ws.Range("R3").Value = 14316.26709487648
$ws.Range("S3").Value = 0.4926993874977079
$ws.Range("T3").Value = 0.4926993874977079
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.04537966666667
$ws.Range("H4").Value = 54.13613900000001
$ws.Range("I4").Value = 0.6797959733292525
$ws.Range("J4").Value = 0.6797959733292525
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.24063
$ws.Range("N4").Value = 0.72189
$ws.Range("O4").Value = 0.001978481285600361
$ws.Range("P4").Value = 0.001978481285600361
$ws.Range("Q4").Value = 4.34225970919
$ws.Range("R4").Value = 39.08033738271001
$ws.Range("S4").Value = 0.001344963611258408
$ws.Range("T4").Value = 0.001344963611258408
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 18.04537966666667
$ws.Range("H5").Value = 54.13613900000001
$ws.Range("I5").Value = 0.6797959733292525
$ws.Range("J5").Value = 0.6797959733292525
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.53182233333333
$ws.Range("N5").Value = 91.595467
$ws.Range("O5").Value = 0.2510353617660938
$ws.Range("P5").Value = 0.2510353617660938
$ws.Range("Q5").Value = 550.9583259202127
$ws.Range("R5").Value = 4958.624933281913
$ws.Range("S5").Value = 0.1706528280918428
$ws.Range("T5").Value = 0.1706528280918428
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.6001993333333334
$ws.Range("H6").Value = 1.800598
$ws.Range("I6").Value = 0.02261039099934159
$ws.Range("J6").Value = 0.02261039099934159
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.701354
$ws.Range("N6").Value = 8.104061999999999
$ws.Range("O6").Value = 0.02221077311549548
$ws.Range("P6").Value = 0.02221077311549548
$ws.Range("Q6").Value = 1.621350869897333
$ws.Range("R6").Value = 14.592157829076
$ws.Range("S6").Value = 0.0005021942645390173
$ws.Range("T6").Value = 0.0005021942645390173
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.6001993333333334
$ws.Range("H7").Value = 1.800598
$ws.Range("I7").Value = 0.02261039099934159
$ws.Range("J7").Value = 0.02261039099934159
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("N7").Value = 264.449356
$ws.Range("O7").Value = 0.7247753838328104
$ws.Range("P7").Value = 0.7247753838328105
$ws.Range("Q7").Value = 52.90744239054312
$ws.Range("R7").Value = 476.166981514888
$ws.Range("S7").Value = 0.01638745481515773
$ws.Range("T7").Value = 0.01638745481515773
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.6001993333333334
$ws.Range("H8").Value = 1.800598
$ws.Range("I8").Value = 0.02261039099934159
$ws.Range("J8").Value = 0.02261039099934159
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.24063
$ws.Range("N8").Value = 0.72189
$ws.Range("O8").Value = 0.001978481285600361
$ws.Range("P8").Value = 0.001978481285600361
$ws.Range("Q8").Value = 0.14442596558
$ws.Range("R8").Value = 1.29983369022
$ws.Range("S8").Value = [double]"4.473423545230419E-05"
$ws.Range("T8").Value = [double]"4.473423545230419E-05"
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.6001993333333334
$ws.Range("H9").Value = 1.800598
$ws.Range("I9").Value = 0.02261039099934159
$ws.Range("J9").Value = 0.02261039099934159
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.53182233333333
$ws.Range("N9").Value = 91.595467
$ws.Range("O9").Value = 0.2510353617660938
$ws.Range("P9").Value = 0.2510353617660938
$ws.Range("Q9").Value = 18.32517940991845
$ws.Range("R9").Value = 164.926614689266
$ws.Range("S9").Value = 0.005676007684192548
$ws.Range("T9").Value = 0.005676007684192548
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.690054333333333
$ws.Range("H10").Value = 14.070163
$ws.Range("I10").Value = 0.1766812397072912
$ws.Range("J10").Value = 0.1766812397072912
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.701354
$ws.Range("N10").Value = 8.104061999999999
$ws.Range("O10").Value = 0.02221077311549548
$ws.Range("P10").Value = 0.02221077311549548
$ws.Range("Q10").Value = 12.66949703356733
$ws.Range("R10").Value = 114.025473302106
$ws.Range("S10").Value = 0.003924226928903115
$ws.Range("T10").Value = 0.003924226928903116
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.690054333333333
$ws.Range("H11").Value = 14.070163
$ws.Range("I11").Value = 0.1766812397072912
$ws.Range("J11").Value = 0.1766812397072912
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 88.14978533333333
$ws.Range("N11").Value = 264.449356
$ws.Range("O11").Value = 0.7247753838328104
$ws.Range("P11").Value = 0.7247753838328105
$ws.Range("Q11").Value = 413.4272826850031
$ws.Range("R11").Value = 3720.845544165028
$ws.Range("S11").Value = 0.1280542133249087
$ws.Range("T11").Value = 0.1280542133249088
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.690054333333333
$ws.Range("H12").Value = 14.070163
$ws.Range("I12").Value = 0.1766812397072912
$ws.Range("J12").Value = 0.1766812397072912
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.24063
$ws.Range("N12").Value = 0.72189
$ws.Range("O12").Value = 0.001978481285600361
$ws.Range("P12").Value = 0.001978481285600361
$ws.Range("Q12").Value = 1.12856777423
$ws.Range("R12").Value = 10.15710996807
$ws.Range("S12").Value = 0.000349560526277547
$ws.Range("T12").Value = 0.000349560526277547
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.690054333333333
$ws.Range("H13").Value = 14.070163
$ws.Range("I13").Value = 0.1766812397072912
$ws.Range("J13").Value = 0.1766812397072912
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.53182233333333
$ws.Range("N13").Value = 91.595467
$ws.Range("O13").Value = 0.2510353617660938
$ws.Range("P13").Value = 0.2510353617660938
$ws.Range("Q13").Value = 143.1959056390135
$ws.Range("R13").Value = 1288.763150751121
$ws.Range("S13").Value = 0.04435323892720178
$ws.Range("T13").Value = 0.04435323892720178
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.209654333333333
$ws.Range("H14").Value = 9.628962999999999
$ws.Range("I14").Value = 0.1209123959641148
$ws.Range("J14").Value = 0.1209123959641148
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.701354
$ws.Range("N14").Value = 8.104061999999999
$ws.Range("O14").Value = 0.02221077311549548
$ws.Range("P14").Value = 0.02221077311549548
$ws.Range("Q14").Value = 8.670412571967331
$ws.Range("R14").Value = 78.03371314770598
$ws.Range("S14").Value = 0.002685557793609905
$ws.Range("T14").Value = 0.002685557793609905
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.209654333333333
$ws.Range("H15").Value = 9.628962999999999
$ws.Range("I15").Value = 0.1209123959641148
$ws.Range("J15").Value = 0.1209123959641148
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 88.14978533333333
$ws.Range("N15").Value = 264.449356
$ws.Range("O15").Value = 0.7247753838328104
$ws.Range("P15").Value = 0.7247753838328105
$ws.Range("Q15").Value = 282.9303404775364
$ws.Range("R15").Value = 2546.373064297828
$ws.Range("S15").Value = 0.08763432819503605
$ws.Range("T15").Value = 0.08763432819503607
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.209654333333333
$ws.Range("H16").Value = 9.628962999999999
$ws.Range("I16").Value = 0.1209123959641148
$ws.Range("J16").Value = 0.1209123959641148
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.24063
$ws.Range("N16").Value = 0.72189
$ws.Range("O16").Value = 0.001978481285600361
$ws.Range("P16").Value = 0.001978481285600361
$ws.Range("Q16").Value = 0.7723391222299999
$ws.Range("R16").Value = 6.951052100069999
$ws.Range("S16").Value = 0.0002392229126121018
$ws.Range("T16").Value = 0.0002392229126121018
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.209654333333333
$ws.Range("H17").Value = 9.628962999999999
$ws.Range("I17").Value = 0.1209123959641148
$ws.Range("J17").Value = 0.1209123959641148
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.53182233333333
$ws.Range("N17").Value = 91.595467
$ws.Range("O17").Value = 0.2510353617660938
$ws.Range("P17").Value = 0.2510353617660938
$ws.Range("Q17").Value = 97.99659585674677
$ws.Range("R17").Value = 881.9693627107209
$ws.Range("S17").Value = 0.03035328706285674
$ws.Range("T17").Value = 0.03035328706285674